$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.290.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.577.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.57%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.65%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.82%  "

$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("E10").Value = "  +1.72%  "

$ws.Range("E11").Value = "  -0.78%  "

$ws.Range("E12").Value = "  -0.66%  "

$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.040.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.053.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.24%  "

$ws.Range("E16").Value = "  +5.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.582.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.26%  "

$ws.Range("E18").Value = "  +2.86%  "

$ws.Range("E19").Value = "  +2.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("E21").Value = "  -1.46%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.16%  "

$ws.Range("E24").Value = "  +0.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.33%  "

$ws.Range("E26").Value = "  -1.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "556.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "

$ws.Range("E28").Value = "  -1.70%  "

$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  -1.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0855"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("E33").Value = "  -1.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "166.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.411"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.39%  "

$ws.Range("E39").Value = "  -2.44%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "165.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.50%  "

$ws.Range("E43").Value = "  +4.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.38%  "

$ws.Range("E45").Value = "  +1.85%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.63%  "

$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("E48").Value = "  +1.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0960"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("E51").Value = "  +19.26%  "
